$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 630
$ws1.Range("F3").Value = 578
$ws1.Range("F5").Value = 32
$ws1.Range("F6").Value = 127
$ws1.Range("F10").Value = 5079
$ws1.Range("F11").Value = 4729
$ws1.Range("F13").Value = 32
$ws1.Range("F15").Value = 45
$ws1.Range("F16").Value = 179

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 73

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 630
$ws4.Range("F3").Value = 578
$ws4.Range("F5").Value = 32
$ws4.Range("F6").Value = 127
$ws4.Range("F10").Value = 5079
$ws4.Range("F11").Value = 4729
$ws4.Range("F13").Value = 32
$ws4.Range("F15").Value = 45
$ws4.Range("F16").Value = 179
$ws4.Range("F17").Value = 73
